# Aggiornamento fino al 27/05 - append new daily rows (256-269) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row, date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(256, 44330, 2, 5, 46.06172270842929),
    @(257, 44331, 0, 5, 46.06172270842929),
    @(258, 44332, 0, 3, 27.63703362505758),
    @(259, 44333, 1, 4, 36.84937816674343),
    @(260, 44334, 0, 4, 36.84937816674343),
    @(261, 44335, 1, 5, 46.06172270842929),
    @(262, 44336, 0, 4, 36.84937816674343),
    @(263, 44337, 0, 2, 18.42468908337172),
    @(264, 44338, 2, 4, 36.84937816674343),
    @(265, 44339, 0, 4, 36.84937816674343),
    @(266, 44340, 1, 4, 36.84937816674343),
    @(267, 44341, 0, 4, 36.84937816674343),
    @(268, 44342, 0, 3, 27.63703362505758),
    @(269, 44343, 0, 3, 27.63703362505758)
)

# Copy the style/format from the last existing row (255) in column A so that
# the new date cells keep the same number format / font / border / alignment.
$ws.Range("A255").Copy() | Out-Null

foreach ($entry in $data) {
    $r = $entry[0]

    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("A$r").Value = $entry[1]
    $ws.Range("B$r").Value = $entry[2]
    $ws.Range("C$r").Value = $entry[3]
    $ws.Range("D$r").Value = $entry[4]
}

$excel.CutCopyMode = 0
